$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 452) { $lastRow = 452 }

$ws.Range("C2:C$lastRow").Value = 45177
